# edit.ps1
# Applies the "Updated symbol list" crypto price/volume refresh commit.
# Column D = Price, Column E = Volume(1h). Both are stored as text in the
# worksheet, so we force each cell's NumberFormat to Text ("@") before
# writing its value, preventing Excel from auto-converting the strings to
# numbers (which would strip the "%" suffix / change numeric precision).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.52%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.46%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.262"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.56%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08087"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.507"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.77%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.646"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.98%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.929"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.46%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9360"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.19%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1349"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "24.27%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1975"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.91%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09173"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.57%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03498"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "4.80%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09586"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.00%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001357"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.92%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006430"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "10.78%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.365"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.92%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.09%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.038"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "12.86%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.52%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "10.91%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.05%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001222"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.99%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004317"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.44%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001191"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.41%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003991"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02486"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05200"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.76%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007723"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.89%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1428"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.80%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009224"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.53%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002173"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01053"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "32.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006638"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.64%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003340"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.82%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "148.11%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
